$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows that previously held "RM 232" and "SC 5" are removed; every
# row below them shifts up by two, which is how "RM 232" disappears and the
# former "SC 5" row becomes the new row 26, etc. Deleting the rows keeps all
# the still-correct values (for rows that don't change) intact automatically.
$ws.Rows("26:27").Delete()

# --- Value corrections (post row-shift) -------------------------------------

# F2: value cleared (now blank/missing)
$ws.Range("F2").Value = ""

# F5: newly filled in
$ws.Range("F5").Value = 17.66

# D6, F6: newly filled in
$ws.Range("D6").Value = -14.2
$ws.Range("F6").Value = 16.43

# D8: cleared
$ws.Range("D8").Value = ""

# F9: cleared
$ws.Range("F9").Value = ""

# F10: cleared
$ws.Range("F10").Value = ""

# D12: newly filled in
$ws.Range("D12").Value = -14.1

# D14: cleared
$ws.Range("D14").Value = ""

# D17: newly filled in
$ws.Range("D17").Value = -14.7

# D18: newly filled in
$ws.Range("D18").Value = -15.2

# D19: cleared
$ws.Range("D19").Value = ""

# D20: cleared
$ws.Range("D20").Value = ""

# D23: newly filled in
$ws.Range("D23").Value = -13.9

# F24: newly filled in
$ws.Range("F24").Value = 16.78

# Row 26 (now "SC 92" after the shift) becomes "SC 5" with its own values
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

# Row 27 (now "SC 101" after the shift): C filled in, D cleared
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = ""

# Row 28 (now "SC 105" after the shift): C and F cleared
$ws.Range("C28").Value = ""
$ws.Range("F28").Value = ""

# Row 29 (now "SC 119" after the shift): C cleared
$ws.Range("C29").Value = ""

# Row 30 (now "SC 120" after the shift): C and F filled in
$ws.Range("C30").Value = 11.4
$ws.Range("F30").Value = 16.89

# Row 32 (now "SC 193" after the shift): C cleared
$ws.Range("C32").Value = ""
